$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-19: update D, E, F, G, H values ---
$rows2to19 = @(
    @(2, 0.1142586985524135, 0.1065029325318719, 0.08528972426563008, 0.1105213794923197, 0.1121894524774864),
    @(3, 0.1187580543085927, 0.2141657941562525, 0.1673868894695735, 0.256203685924516, 0.2457789888377989),
    @(4, 0.1725693685531227, 0.3419621502174148, 0.3269998128819627, 0.4845143733610661, 0.5543580886384267),
    @(5, 0.1990363627893637, 0.6133483807311121, 0.7165058510440035, 1.377173441754483, 1.191427732428485),
    @(6, 0.2060247184665036, 1.06290747902815, 1.259567647197757, 1.996983857992604, 2.865244898738017),
    @(7, 0.3020367966961433, 2.403847289496705, 1.947104405873099, 6.336478816774412, 5.859442856540852),
    @(8, 0.3440681853332487, 3.770094562520921, 3.774921837413167, 13.88154630555728, 12.78056414330938),
    @(9, 0.532300744986195, 6.143512738541525, 6.952590901171519, 25.13311559123643, 25.79522056605397),
    @(10, 0.5385402165298068, 11.8960194251739, 10.20541107356204, 42.27788787617155, 39.88653473205823),
    @(11, 0.6510933183425628, 20.26049679386849, 22.65766035830179, 49.99008797771357, 54.54728033762993),
    @(12, 0.665021419322003, 30.03083656705767, 27.84439581032989, 71.48444916046137, 63.8137482365326),
    @(13, 0.8978264846541173, 44.27852814319333, 39.12678193089796, 74.09910510755705, 87.26450807714986),
    @(14, 1.033283150342568, 45.67080936938287, 59.907607919771, 90.78618417746785, 89.54929501879141),
    @(15, 1.263421402216381, 70.73665541828065, 60.78965032427504, 87.98178657913039, 91.06999213085005),
    @(16, 1.392895266056499, 73.07330058366908, 65.84091318132573, 80.99310163724661, 96.87109819069866),
    @(17, 1.874790045537552, 89.65104007731702, 78.2622059198177, 78.30611838299704, 90.61741422682927),
    @(18, 2.582406295991133, 75.15774815778268, 81.72342201720647, 80.98696124132879, 80.6606642161201),
    @(19, 2.835766969419082, 85.9569115277357, 98.33439195760907, 89.72866519261916, 85.73051280185452)
)
foreach ($row in $rows2to19) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 5).Value = $row[2]
    $ws.Cells.Item($r, 6).Value = $row[3]
    $ws.Cells.Item($r, 7).Value = $row[4]
    $ws.Cells.Item($r, 8).Value = $row[5]
}

# --- Rows 20-25: clear G/H, update D, add E, F ---
$ws.Range("G20:H25").ClearContents()
$rows20to25 = @(
    @(20, 3.273173129454402, 87.34119705876086, 84.01270653704678),
    @(21, 3.412741350136176, 89.46740484542516, 90.03151993377432),
    @(22, 5.3852799659899, 88.17685540040674, 80.74747839747401),
    @(23, 6.042979635111466, 96.14246368279656, 97.35017504148126),
    @(24, 5.900160113037796, 82.16205014886536, 87.04332756701884),
    @(25, 7.761892011177026, 106.7506321427255, 79.06541391240165)
)
foreach ($row in $rows20to25) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 5).Value = $row[2]
    $ws.Cells.Item($r, 6).Value = $row[3]
}

# --- Rows 26-39: update D only ---
$rows26to39 = @(
    @(26, 10.59995211647931),
    @(27, 12.45388250053755),
    @(28, 13.73740824591042),
    @(29, 13.5592138042465),
    @(30, 18.03862714393862),
    @(31, 22.90648708505597),
    @(32, 27.42614979530502),
    @(33, 26.46724133078785),
    @(34, 34.61380459902725),
    @(35, 40.25427505071339),
    @(36, 38.98230435045711),
    @(37, 43.31366938449019),
    @(38, 52.98292123012438),
    @(39, 48.38875170931617)
)
foreach ($row in $rows26to39) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
}

# --- Rows 40-71: new rows with A, B, D ---
$newRows = @(
    @(40, 38, 38, 58.33749036609034),
    @(41, 39, 39, 68.90245716432005),
    @(42, 40, 40, 71.97372000110504),
    @(43, 41, 41, 71.5657339775115),
    @(44, 42, 42, 64.61622218675642),
    @(45, 43, 43, 66.11605913247193),
    @(46, 44, 44, 74.78083973404182),
    @(47, 45, 45, 71.26286997381952),
    @(48, 46, 46, 74.63562464915874),
    @(49, 47, 47, 86.62614569471329),
    @(50, 48, 48, 86.06058270085954),
    @(51, 49, 49, 80.24467644362585),
    @(52, 50, 50, 81.92397510831036),
    @(53, 51, 51, 86.34340461651389),
    @(54, 52, 52, 74.55565215089713),
    @(55, 53, 53, 86.73241101452844),
    @(56, 54, 54, 76.69989009450561),
    @(57, 55, 55, 96.84993585147643),
    @(58, 56, 56, 86.80312593206111),
    @(59, 57, 57, 90.37439563609601),
    @(60, 58, 58, 85.79923199034162),
    @(61, 59, 59, 82.85300710027032),
    @(62, 60, 60, 89.38287341712579),
    @(63, 61, 61, 91.11236394190499),
    @(64, 62, 62, 99.17707865924018),
    @(65, 63, 63, 95.45427816986368),
    @(66, 64, 64, 83.66023796640106),
    @(67, 65, 65, 89.73924185571688),
    @(68, 66, 66, 100.5520129959937),
    @(69, 67, 67, 93.44211936242009),
    @(70, 68, 68, 83.92095380673342),
    @(71, 69, 69, 74.67877594807399)
)
foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Apply the same style as A39 (bold, bordered, centered) to new A column cells
$ws.Range("A39").Copy()
$ws.Range("A40:A71").PasteSpecial(-4122)
$excel.CutCopyMode = $false
